$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2351421188630491
$ws.Range("C2").Value = 0.4935400516795866
$ws.Range("J2").Value = 0.02325581395348837
$ws.Range("P2").Value = 0.1679586563307494
$ws.Range("S2").Value = 0.08010335917312661
$ws.Range("B3").Value = 0.01492537313432836
$ws.Range("C3").Value = 0.02487562189054726
$ws.Range("J3").Value = 0.04477611940298507
$ws.Range("P3").Value = 0.7263681592039801
$ws.Range("S3").Value = 0.1890547263681592
$ws.Range("J4").Value = 0.08771929824561403
$ws.Range("O4").Value = 0.01754385964912281
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2280701754385965
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.07373271889400922
$ws.Range("D6").Value = 0.009216589861751152
$ws.Range("F6").Value = 0.05529953917050692
$ws.Range("J6").Value = 0.2350230414746544
$ws.Range("O6").Value = 0.02304147465437788
$ws.Range("Q6").Value = 0.1751152073732719
$ws.Range("R6").Value = 0.06912442396313365
$ws.Range("S6").Value = 0.3594470046082949
$ws.Range("B7").Value = 0.1123595505617977
$ws.Range("D7").Value = 0.01123595505617977
$ws.Range("F7").Value = 0.06741573033707865
$ws.Range("J7").Value = 0.1348314606741573
$ws.Range("O7").Value = 0.03932584269662921
$ws.Range("Q7").Value = 0.2078651685393259
$ws.Range("R7").Value = 0.06741573033707865
$ws.Range("S7").Value = 0.3595505617977528
$ws.Range("B8").Value = 0.09657320872274143
$ws.Range("D8").Value = 0.01557632398753894
$ws.Range("E8").Value = 0.003115264797507788
$ws.Range("F8").Value = 0.0529595015576324
$ws.Range("J8").Value = 0.1526479750778816
$ws.Range("O8").Value = 0.02180685358255452
$ws.Range("Q8").Value = 0.2305295950155763
$ws.Range("R8").Value = 0.1059190031152648
$ws.Range("S8").Value = 0.3208722741433022
$ws.Range("B9").Value = 0.1541666666666667
$ws.Range("D9").Value = 0.008333333333333333
$ws.Range("E9").Value = 0.004166666666666667
$ws.Range("F9").Value = 0.05
$ws.Range("J9").Value = 0.1208333333333333
$ws.Range("O9").Value = 0.04166666666666666
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.07083333333333333
$ws.Range("S9").Value = 0.3833333333333334
$ws.Range("B10").Value = 0.1216828478964401
$ws.Range("D10").Value = 0.03042071197411003
$ws.Range("E10").Value = 0.001941747572815534
$ws.Range("F10").Value = 0.06666666666666667
$ws.Range("J10").Value = 0.112621359223301
$ws.Range("O10").Value = 0.01812297734627832
$ws.Range("Q10").Value = 0.2394822006472492
$ws.Range("R10").Value = 0.09061488673139159
$ws.Range("S10").Value = 0.3184466019417476
$ws.Range("G11").Value = 0.1176470588235294
$ws.Range("J11").Value = 0.1245674740484429
$ws.Range("K11").Value = 0.1799307958477509
$ws.Range("L11").Value = 0.5709342560553633
$ws.Range("S11").Value = 0.006920415224913495
$ws.Range("G12").Value = 0.7183908045977011
$ws.Range("J12").Value = 0.1896551724137931
$ws.Range("K12").Value = 0.01724137931034483
$ws.Range("L12").Value = 0.05172413793103448
$ws.Range("S12").Value = 0.02298850574712644
$ws.Range("G13").Value = 0.6216216216216216
$ws.Range("J13").Value = 0.3783783783783784
$ws.Range("F15").Value = 0.02097902097902098
$ws.Range("H15").Value = 0.1118881118881119
$ws.Range("I15").Value = 0.06993006993006994
$ws.Range("J15").Value = 0.3601398601398602
$ws.Range("K15").Value = 0.05944055944055944
$ws.Range("M15").Value = 0.003496503496503497
$ws.Range("N15").Value = 0.003496503496503497
$ws.Range("O15").Value = 0.05244755244755245
$ws.Range("S15").Value = 0.3181818181818182
$ws.Range("F16").Value = 0.01244813278008299
$ws.Range("H16").Value = 0.1161825726141079
$ws.Range("I16").Value = 0.0954356846473029
$ws.Range("J16").Value = 0.4813278008298755
$ws.Range("K16").Value = 0.1120331950207469
$ws.Range("M16").Value = 0.008298755186721992
$ws.Range("O16").Value = 0.03734439834024896
$ws.Range("S16").Value = 0.1369294605809129
$ws.Range("F17").Value = 0.00544464609800363
$ws.Range("H17").Value = 0.1016333938294011
$ws.Range("I17").Value = 0.1288566243194192
$ws.Range("J17").Value = 0.5099818511796733
$ws.Range("K17").Value = 0.06896551724137931
$ws.Range("M17").Value = 0.01270417422867514
$ws.Range("O17").Value = 0.07259528130671507
$ws.Range("S17").Value = 0.0998185117967332
$ws.Range("F18").Value = 0.01388888888888889
$ws.Range("H18").Value = 0.08796296296296297
$ws.Range("I18").Value = 0.1111111111111111
$ws.Range("J18").Value = 0.4166666666666667
$ws.Range("K18").Value = 0.1064814814814815
$ws.Range("M18").Value = 0.02314814814814815
$ws.Range("N18").Value = 0.004629629629629629
$ws.Range("O18").Value = 0.09722222222222222
$ws.Range("S18").Value = 0.1388888888888889
$ws.Range("F19").Value = 0.01115537848605578
$ws.Range("H19").Value = 0.149003984063745
$ws.Range("I19").Value = 0.0796812749003984
$ws.Range("J19").Value = 0.4254980079681275
$ws.Range("K19").Value = 0.099601593625498
$ws.Range("M19").Value = 0.01832669322709163
$ws.Range("O19").Value = 0.09322709163346614
$ws.Range("S19").Value = 0.1235059760956175
